$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume figures (GitHub Actions scheduled refresh)
# Some "Price" values are purely numeric-looking strings (e.g. "557.79");
# force them to remain plain text (matching the original inline-string cells)
# instead of being auto-converted to numbers, then restore the default style.

$ws.Range('D2').Value = '59.826.73'
$ws.Range('E2').Value = '  -0.05%  '
$ws.Range('D3').Value = '2.369.39'
$ws.Range('E3').Value = '  -1.67%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '557.79'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '133.26'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.84%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -1.45%  '
$ws.Range('E9').Value = '  -0.30%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.63'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.80%  '
$ws.Range('E11').Value = '  +1.13%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.342'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.34%  '
$ws.Range('E13').Value = '  -3.89%  '
$ws.Range('D14').Value = '2.790.92'
$ws.Range('E14').Value = '  -1.76%  '
$ws.Range('D15').Value = '59.761.74'
$ws.Range('E16').Value = '  -0.03%  '
$ws.Range('D17').Value = '2.367.92'
$ws.Range('E17').Value = '  -1.76%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.07'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.16%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '320.25'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.56%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.65'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.33%  '
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '64.31'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.53%  '
$ws.Range('E24').Value = '  -0.48%  '
$ws.Range('E25').Value = '  -0.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.40'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.69%  '
$ws.Range('E27').Value = '  -0.50%  '
$ws.Range('E28').Value = '  +2.00%  '
$ws.Range('D29').Value = '0.0₃0757'
$ws.Range('E29').Value = '  -1.76%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '170.34'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.63%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.08'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.87%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.12'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +10.53%  '
$ws.Range('E33').Value = '  -1.02%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.11'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.72%  '
$ws.Range('E36').Value = '  +0.68%  '
$ws.Range('E37').Value = '  +0.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.13'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.22%  '
$ws.Range('E39').Value = '  -1.17%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '318.86'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.76%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '38.60'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.20%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '144.50'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.17%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.53'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.40%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0966'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.23%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '19.43'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.42%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0510'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.77%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.568'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.75%  '
$ws.Range('E48').Value = '  -2.75%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '11.06'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.15%  '
$ws.Range('E50').Value = '  -0.35%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.54'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.58%  '
